## QuyenTX: Finished W03, verifying W04
## Applies the "ServiceImporting" verification-block edits (W04) plus
## moving the active/selected tab from Login -> ServiceImporting.

$wb = $excel.ActiveWorkbook

$wsLogin = $wb.Worksheets.Item("Login")
$wsSvc   = $wb.Worksheets.Item("ServiceImporting")

# ---------------------------------------------------------------------
# 1. ServiceImporting (sheet6) content changes
# ---------------------------------------------------------------------

# Column C needs to widen to fit the new "Project Key" / "SP-1" column.
$wsSvc.Columns.Item(3).ColumnWidth = 13

# Two blank spacer rows (4 & 5) above the existing W03 block.
$wsSvc.Rows.Item(4).RowHeight = $wsSvc.Rows.Item(4).RowHeight
$wsSvc.Rows.Item(5).RowHeight = $wsSvc.Rows.Item(5).RowHeight

# Existing W03 block (rows 6-8) gains a "Project Key" column (C/D).
$wsSvc.Range("C6").Value = "Project Key"
$c6 = $wsSvc.Range("C6")
$a6 = $wsSvc.Range("A6")
$a6.Copy() | Out-Null
$c6.PasteSpecial(-4122) | Out-Null   # xlPasteFormats (A6 is the yellow header style)
$b6 = $wsSvc.Range("B6")
$b6.Copy() | Out-Null
$c6.PasteSpecial(-4122) | Out-Null   # xlPasteFormats (B6's orange style is what C6 actually needs)
$wsSvc.Range("C6").Value = "Project Key"

$d6 = $wsSvc.Range("D6")
$c7Old = $wsSvc.Range("C7")
$c7Old.Copy() | Out-Null
$d6.PasteSpecial(-4122) | Out-Null

$wsSvc.Range("B7").Value = "Selenium Practice"
$wsSvc.Range("C7").Value = "SP-1"
$wsSvc.Rows.Item(7).RowHeight = 30

$d7 = $wsSvc.Range("D7")
$c7Old.Copy() | Out-Null
$d7.PasteSpecial(-4122) | Out-Null

# Move the "W03" marker cell from C8 to D8.
$c8 = $wsSvc.Range("C8")
$d8 = $wsSvc.Range("D8")
$c8.Cut($d8) | Out-Null
$c7Old.Copy() | Out-Null
$c8.PasteSpecial(-4122) | Out-Null

# New W04 verification block (rows 10-12), mirroring the W03 block above.
$wsSvc.Range("A10").Value = "W02"
$a10 = $wsSvc.Range("A10")
$a6.Copy() | Out-Null
$a10.PasteSpecial(-4122) | Out-Null
$wsSvc.Range("A10").Value = "W02"

$wsSvc.Range("B10").Value = "Project key"
$b10 = $wsSvc.Range("B10")
$b6.Copy() | Out-Null
$b10.PasteSpecial(-4122) | Out-Null
$wsSvc.Range("B10").Value = "Project key"

$wsSvc.Range("C10").Value = "Project Name"
$c10 = $wsSvc.Range("C10")
$b6.Copy() | Out-Null
$c10.PasteSpecial(-4122) | Out-Null
$wsSvc.Range("C10").Value = "Project Name"

$wsSvc.Range("D10").Value = "Project Department"
$d10 = $wsSvc.Range("D10")
$b6.Copy() | Out-Null
$d10.PasteSpecial(-4122) | Out-Null
$wsSvc.Range("D10").Value = "Project Department"

$wsSvc.Range("E10").Value = "Project Type"
$e10 = $wsSvc.Range("E10")
$b6.Copy() | Out-Null
$e10.PasteSpecial(-4122) | Out-Null
$wsSvc.Range("E10").Value = "Project Type"

$f10 = $wsSvc.Range("F10")
$c7Old.Copy() | Out-Null
$f10.PasteSpecial(-4122) | Out-Null

$wsSvc.Rows.Item(10).RowHeight = 30

$wsSvc.Range("B11").Value = "SP"
$wsSvc.Range("C11").Value = "Selenium Practice"
$wsSvc.Range("D11").Value = "VSII / TSC"
$wsSvc.Range("E11").Value = "JIRA"
$wsSvc.Rows.Item(11).RowHeight = 30

foreach ($addr in @("A11","B11","C11","D11","E11","F11")) {
    $cell = $wsSvc.Range($addr)
    $c7Old.Copy() | Out-Null
    $cell.PasteSpecial(-4122) | Out-Null
}
$wsSvc.Range("B11").Value = "SP"
$wsSvc.Range("C11").Value = "Selenium Practice"
$wsSvc.Range("D11").Value = "VSII / TSC"
$wsSvc.Range("E11").Value = "JIRA"

$f12 = $wsSvc.Range("F12")
$f12.Value = "W02"
$a6.Copy() | Out-Null
$f12.PasteSpecial(-4122) | Out-Null
$f12.Value = "W02"

foreach ($addr in @("A12","B12","C12","D12","E12")) {
    $cell = $wsSvc.Range($addr)
    $c7Old.Copy() | Out-Null
    $cell.PasteSpecial(-4122) | Out-Null
}

$wsSvc.PageSetup.Orientation = 1   # xlPortrait

# ---------------------------------------------------------------------
# 2. Selection / active-tab bookkeeping
# ---------------------------------------------------------------------
$wsLogin.Range("D15").Select() | Out-Null
$wsSvc.Activate() | Out-Null
$wsSvc.Range("G7").Select() | Out-Null
